# Lesson 6.2 Trees - update cached date-field text from 10/13/2015 to
# 10/19/2015 across the notes master, slide master and all slide layouts,
# and rewrite the "Next Steps" slide's first bullet to add the file
# reference, splitting it into two paragraphs.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "10/19/2015"

# Notes master "Date Placeholder" (datetimeFigureOut field).
Set-DatePlaceholderText $p.NotesMaster.Shapes $newDate

# Slide master "Date Placeholder" (datetime1 field).
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout's own "Date Placeholder" (datetime1 field).
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 12 ("Next Steps") - Content Placeholder 2: replace the single
# "If you have questions..." paragraph with two paragraphs, adding a
# reference to the Examples folder file before it.
$slide12 = $p.Slides.Item(12)
$contentShape = $null
for ($i = 1; $i -le $slide12.Shapes.Count; $i++) {
    $shp = $slide12.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 2") {
        $contentShape = $shp
    }
}

$tr = $contentShape.TextFrame.TextRange
$tr.Paragraphs(1, 1).Text = "Study the file 06-2-trees.rkt in the Examples folder.`rIf you have questions about this lesson, ask them on the Discussion Board"
